# Update "nombre_aides" (column C) and "montant_total" (column D) figures
# for the 2020-12-30 Fonds de solidarite volet 2 regional/NAF extract.
#
# Both columns are stored as text in the workbook (e.g. "130281.00" keeps
# its trailing zeroes), so each cell is forced to Text format before the
# new value is written and the style is reset to "Normal" right after so
# that no stray formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> new nombre_aides (col C), new montant_total (col D)
$updates = @(
    @{ Row = 21;  C = "48";   D = "130281.00" }
    @{ Row = 24;  C = "494";  D = "3806709.26" }
    @{ Row = 39;  C = "578";  D = "6190211.97" }
    @{ Row = 47;  C = "96";   D = "1028987.72" }
    @{ Row = 130; C = "1125"; D = "9105363.67" }
    @{ Row = 214; C = "1006"; D = "9724988.80" }
    @{ Row = 222; C = "208";  D = "2189031.09" }
    @{ Row = 245; C = "475";  D = "3465751.80" }
    @{ Row = 249; C = "104";  D = "809165.53" }
    @{ Row = 253; C = "92";   D = "818127.57" }
    @{ Row = 261; C = "1710"; D = "11976409.68" }
)

foreach ($u in $updates) {
    $cellC = $ws.Cells.Item($u.Row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $u.C
    $cellC.Style = "Normal"

    $cellD = $ws.Cells.Item($u.Row, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $u.D
    $cellD.Style = "Normal"
}
